$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" column (C) for data rows 2 through 39
# from serial date 45171 (2023-09-02) to 45172 (2023-09-03)
for ($row = 2; $row -le 39; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = 45172
}
